$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.0196232996448211, 0.7935034002069985, 1.1052023992871707, 2.941813241995588)"
$ws.Range("C2").Value = "NIG(1.4037115723397222, 1.0644434441592936, 4.103300018528145, 6.223297131827735)"
$ws.Range("D2").Value = "NCT(3.437347556551492, 1.3879533321118818, -0.12946529077385394, 2.716830113383709)"
$ws.Range("E2").Value = "NIG(1.5171266651900308, 1.223947504252327, 3.247538070580579, 5.680925933650009)"
